$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Round the Ost (Q) and Nord (R) coordinate values to whole numbers for rows 12-14
$ws.Range("Q12").Value = 374850
$ws.Range("R12").Value = 6871061

$ws.Range("Q13").Value = 375047
$ws.Range("R13").Value = 6871264

$ws.Range("Q14").Value = 374954
$ws.Range("R14").Value = 6870892

# Remove the Starttid (Z) and Sluttid (AB) values for rows 12-14
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()

$ws.Range("Z13").ClearContents()
$ws.Range("AB13").ClearContents()

$ws.Range("Z14").ClearContents()
$ws.Range("AB14").ClearContents()
